$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 4.c.1 Indicator text updated
$ws.Range("B4").Value = "4.c.1 Proportion of teachers with the minimum required qualifications, by education level"

# Organization website updated
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B10").Font.Name = "Calibri"
$ws.Range("B10").Font.Size = 11

# Reflect the selection left by the editor
$ws.Range("B10").Select()
